# Hjemme passive updated meanEMG legmaxROM
# Updates the "15/16 repetition" columns (B:E) on Ark1 for the
# header row and the two data rows (CON / STR), and moves the
# active selection from B1:AY3 down to B1:E3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# --- Row 1 (header / rep numbers) ---
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# --- Row 2 (CON) ---
$ws.Range("B2").Value = 3.9638822351134952
$ws.Range("C2").Value = 9.4573355776025316
$ws.Range("D2").Value = 11.198267481312284
$ws.Range("E2").Value = 9.8323984339473682

# --- Row 3 (STR) ---
$ws.Range("B3").Value = 3.108782405754853
$ws.Range("C3").Value = 5.7132292175159396
$ws.Range("D3").Value = 15.002643271719867
$ws.Range("E3").Value = 6.6232210289869027

# Narrow the selection to the edited block, matching the updated
# workbook's saved view state (B1:AY3 -> B1:E3).
$excel.Goto($ws.Range("B1:E3"))
